# PowerShell COM-interop script implementing the commit:
# "Textual explaination about strategy and iteratoe patterns added"
#
# The paragraph indices below were derived by simulating the sequence of
# structural edits (inserts/deletes) against the original 13-paragraph body,
# so each $d.Paragraphs(N) below refers to the correct paragraph at the
# moment that particular statement executes.

$d = $word.ActiveDocument

function Set-ParaText($para, [string]$text) {
    $r = $d.Range($para.Range.Start, $para.Range.End - 1)
    $r.Text = $text
}

# 1) P13 "...cities data-base...singleton..." -> emptied
Set-ParaText $d.Paragraphs(13) ''

# 1b) append a brand-new empty (ind=1080) paragraph right after it
$d.Paragraphs(13).Range.InsertParagraphAfter()

# 2) P12 "Singleton - ..." -> "Decorator - "
Set-ParaText $d.Paragraphs(12) 'Decorator - '

# 3) insert a fresh, clean (pStyle a3 / bidi 0 / ind=1080, no rtl) empty
#    paragraph right after P9 ("At this point...") -- this will become the
#    replacement for the old rtl-flagged empty paragraph (old P11).
$d.Paragraphs(9).Range.InsertParagraphAfter()

# 4) delete the old rtl-flagged empty paragraph (now shifted to index 12)
$d.Paragraphs(12).Range.Delete()

# 5) delete "Also, the composer may be altered..." paragraph (now at index 11)
$d.Paragraphs(11).Range.Delete()

# 6) P9 "At this point it is worth considering..." -> "Nevertheless the entity..."
Set-ParaText $d.Paragraphs(9) 'Nevertheless the entity who uses the CsvParser knows it as an interface IParser rather than a class. This makes the format switching t in the future highly convenient.'

# 7) P8 "MusicAverage(SoulMusic , Jazz) = Funk." -> CsvParser sentence
Set-ParaText $d.Paragraphs(8) 'Therefore the parsing mechanism is aggregated in a separate class – CsvParser. '

# 8) P7 (long "UserAverageableDetails" builder detail) -> parsing-format text
Set-ParaText $d.Paragraphs(7) 'In the app it was implemented in order to regulate parsing format. As mentioned in the iterator design pattern, the app used a csv format. Thus the parsing logic was implemented to suit those needs. Nonetheless changing the format to a JSON or XML is certainly a feasible situation as the system grow.'

# 9) P6 "Builder - ..." -> "Strategy - ..."
Set-ParaText $d.Paragraphs(6) 'Strategy – The strategy design pattern is used in order to modify a specific selected part in an overall none-modifiable component. As such it allow to reuse a specific logic in different context. The modifiable part is stored in a different class or in a different function pointer (if supported by the language) and can be swapped upon need or even in runtime.'

# 10) delete P5 (the second, now-redundant empty paragraph)
$d.Paragraphs(5).Range.Delete()

# 11) P4 (empty paragraph) -> gains ind=1080 (54pt) left indent
$d.Paragraphs(4).Range.ParagraphFormat.LeftIndent = 54

# 12) P3 (adapter detail text) -> iterator detail text
Set-ParaText $d.Paragraphs(3) 'In the current project a simple txt file (written in csv format) was used as the cities'' data base. As such it is highly recommended to supply an iterator in order to scan the data-base (especially since the only use of it is for searching a specific city).This the scanning logic is aggregated in one component and the city class is decoupled from it. Also, performance is dramatically improved since searching in a large is no longer requiring the creation of all cities. This structure institutes an infrastructure for the second pattern – strategy. '

# 13) P2 "Adapter - ..." -> "Iterator - ..."
Set-ParaText $d.Paragraphs(2) 'Iterator – The iterator design pattern is used in order to promote reusability, extensibility, maintainability and modularity of iterating trough a specific collection. It does so by defining a scanning algorithm rather than supplying the actual data structure.'

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
